# TC47_Canine_Filter_Breed-YorkshireTerr.xlsx - "Fixed StudyComb for Faceted Filters ICDC"
#
# Replaces the "cartQuery" (column D) Neo4j/Cypher query used by the
# CasesTab / SamplesTab / FilesTab rows with the new, simplified
# "study combination" count query, and refreshes a handful of cosmetic
# view properties (zoom, selection, row heights, column widths) left
# behind by the author's Excel session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core content change -------------------------------------------------
# D2:D4 all shared the same Cypher query text ("cartQuery"); swap it for
# the new StudyComb query.
$newQuery = @"
MATCH (demo:demographic)
WHERE demo.breed IN ['Yorkshire Terrier']
MATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)
OPTIONAL MATCH (c)<-[*]-(samp:sample)
OPTIONAL MATCH (c)<-[*]-(f:file)
RETURN 
	count(DISTINCT(f)) as number_of_files, 
	count(DISTINCT(samp)) as number_of_sample, 
	count(DISTINCT(c)) as number_of_cases, 
	count(DISTINCT(s)) as number_of_study
"@

$ws.Range("D2").Value = $newQuery
$ws.Range("D3").Value = $newQuery
$ws.Range("D4").Value = $newQuery

# --- View / cosmetic changes ---------------------------------------------
# Zoom level moved from 55% to 85%, and the selection moved from C14 to C2.
$excel.ActiveWindow.Zoom = 85
$ws.Range("C2").Select()

# Wrapped rows 2-4 re-autofit to a shorter height now that the query text
# is much shorter (390pt -> 244.8pt).
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Columns picked up tiny width adjustments from the resave (sub-pixel,
# cosmetic only); set the closest values this host's column-width model
# can represent.
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 91.66666666666667
$ws.Columns.Item(3).ColumnWidth = 74.83333333333333
$ws.Columns.Item(4).ColumnWidth = 74.83333333333333
$ws.Columns.Item(5).ColumnWidth = 69.5
$ws.Columns.Item(6).ColumnWidth = 27.666666666666668

Write-Host "Applied StudyComb query update and view refresh."
